$d = $word.ActiveDocument

# The document contains several "<id>...</id>" markers that were each split
# across three separate runs (opening tag / value / closing tag). Re-join
# each of them into a single run (keeping the first run's formatting),
# matching how Word merges runs when you find-and-replace a range that
# spans multiple runs with literal text.

$ids = @("p058v_5", "p059r_1", "p059r_2")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}
